$wb = $excel.ActiveWorkbook
Write-Host "wb type:" $wb.GetType()
$wb | Get-Member | ForEach-Object { Write-Host $_.Name }
